$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates (crypto price/volume list refresh).
# Values that look numeric are apostrophe-prefixed to force text,
# then the style is reset to Normal so no stray number format sticks.

$ws.Range("D2").Value = '45.110.82'
$ws.Range("E2").Value = '  +1.95%  '
$ws.Range("D3").Value = '2.360.44'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'311.81"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = "'107.48"
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").Value = "'8.45"
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("D14").Value = "'0.973"
$ws.Range("E14").Value = '  -4.17%  '
$ws.Range("D15").Value = '2.719.98'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = "'15.21"
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").Value = '2.356.80'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '45.097.91'
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("D19").Value = "'14.21"
$ws.Range("E19").Value = '  +9.40%  '
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("E21").Value = '  -5.75%  '
$ws.Range("D22").Value = "'73.24"
$ws.Range("E22").Value = '  -1.97%  '
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = "'258.92"
$ws.Range("E24").Value = '  -3.73%  '
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = "'11.02"
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("E28").Value = '  -5.42%  '
$ws.Range("D29").Value = "'2.31"
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +6.60%  '
$ws.Range("D31").Value = "'22.31"
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  -6.45%  '
$ws.Range("D33").Value = "'167.46"
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = "'2.99"
$ws.Range("E34").Value = '  +4.85%  '
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").Value = "'3.94"
$ws.Range("E38").Value = '  +4.15%  '
$ws.Range("D39").Value = "'0.0352"
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D42").Value = "'99.83"
$ws.Range("E42").Value = '  -5.00%  '
$ws.Range("D43").Value = "'69.61"
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("E44").Value = '  -5.40%  '
$ws.Range("E45").Value = '  -7.15%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.839.51'
$ws.Range("E47").Value = '  +10.75%  '
$ws.Range("D48").Value = "'83.46"
$ws.Range("E48").Value = '  +7.52%  '
$ws.Range("D49").Value = "'5.69"
$ws.Range("E49").Value = '  +6.92%  '
$ws.Range("D50").Value = "'110.67"
$ws.Range("E50").Value = '  -4.33%  '
$ws.Range("D51").Value = "'9.17"
$ws.Range("E51").Value = '  +1.75%  '

# Reset style on the text-forced cells so no number-format cruft is introduced.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
